$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the selected cell shown in the sheet view (M10 -> O8) ---
$ws.Range("O8").Select()

# --- Remove the AutoFilter from the used range (A1:K112) ---
# Turning AutoFilterMode off drops the worksheet's <autoFilter> element.
$ws.AutoFilterMode = $false

# Deleting the AutoFilter also removes the workbook-level hidden defined
# name Excel maintains for it (_xlnm._FilterDatabase), so the
# <definedNames> block is gone too.
foreach ($n in $wb.Names) {
    $n.Delete()
}

# --- Update column I (row 2 through row 112) ---
# Every row's value moved down by 2. A handful of rows were sitting at the
# 5.01 ceiling (a clamp) in the old data; once shifted down they land below
# the ceiling and reveal the real (previously clamped) values instead of a
# plain "-2" shift, so those rows are written with their literal targets.
$ws.Range("I2").Value = $ws.Range("I2").Value2 - 2
$ws.Range("I3").Value = $ws.Range("I3").Value2 - 2
$ws.Range("I4").Value = $ws.Range("I4").Value2 - 2
$ws.Range("I5").Value = $ws.Range("I5").Value2 - 2
$ws.Range("I6").Value = $ws.Range("I6").Value2 - 2
$ws.Range("I7").Value = $ws.Range("I7").Value2 - 2
$ws.Range("I8").Value = $ws.Range("I8").Value2 - 2
$ws.Range("I9").Value = $ws.Range("I9").Value2 - 2
$ws.Range("I10").Value = $ws.Range("I10").Value2 - 2
$ws.Range("I11").Value = $ws.Range("I11").Value2 - 2
$ws.Range("I12").Value = $ws.Range("I12").Value2 - 2
$ws.Range("I13").Value = 2.71
$ws.Range("I14").Value = $ws.Range("I14").Value2 - 2
$ws.Range("I15").Value = $ws.Range("I15").Value2 - 2
$ws.Range("I16").Value = $ws.Range("I16").Value2 - 2
$ws.Range("I17").Value = $ws.Range("I17").Value2 - 2
$ws.Range("I18").Value = $ws.Range("I18").Value2 - 2
$ws.Range("I19").Value = $ws.Range("I19").Value2 - 2
$ws.Range("I20").Value = $ws.Range("I20").Value2 - 2
$ws.Range("I21").Value = $ws.Range("I21").Value2 - 2
$ws.Range("I22").Value = $ws.Range("I22").Value2 - 2
$ws.Range("I23").Value = $ws.Range("I23").Value2 - 2
$ws.Range("I24").Value = $ws.Range("I24").Value2 - 2
$ws.Range("I25").Value = $ws.Range("I25").Value2 - 2
$ws.Range("I26").Value = $ws.Range("I26").Value2 - 2
$ws.Range("I27").Value = $ws.Range("I27").Value2 - 2
$ws.Range("I28").Value = $ws.Range("I28").Value2 - 2
$ws.Range("I29").Value = $ws.Range("I29").Value2 - 2
$ws.Range("I30").Value = $ws.Range("I30").Value2 - 2
$ws.Range("I31").Value = $ws.Range("I31").Value2 - 2
$ws.Range("I32").Value = $ws.Range("I32").Value2 - 2
$ws.Range("I33").Value = $ws.Range("I33").Value2 - 2
$ws.Range("I34").Value = $ws.Range("I34").Value2 - 2
$ws.Range("I35").Value = $ws.Range("I35").Value2 - 2
$ws.Range("I36").Value = $ws.Range("I36").Value2 - 2
$ws.Range("I37").Value = $ws.Range("I37").Value2 - 2
$ws.Range("I38").Value = $ws.Range("I38").Value2 - 2
$ws.Range("I39").Value = $ws.Range("I39").Value2 - 2
$ws.Range("I40").Value = $ws.Range("I40").Value2 - 2
$ws.Range("I41").Value = $ws.Range("I41").Value2 - 2
$ws.Range("I42").Value = $ws.Range("I42").Value2 - 2
$ws.Range("I43").Value = $ws.Range("I43").Value2 - 2
$ws.Range("I44").Value = $ws.Range("I44").Value2 - 2
$ws.Range("I45").Value = $ws.Range("I45").Value2 - 2
$ws.Range("I46").Value = $ws.Range("I46").Value2 - 2
$ws.Range("I47").Value = $ws.Range("I47").Value2 - 2
$ws.Range("I48").Value = 2.91
$ws.Range("I49").Value = $ws.Range("I49").Value2 - 2
$ws.Range("I50").Value = $ws.Range("I50").Value2 - 2
$ws.Range("I51").Value = $ws.Range("I51").Value2 - 2
$ws.Range("I52").Value = $ws.Range("I52").Value2 - 2
$ws.Range("I53").Value = $ws.Range("I53").Value2 - 2
$ws.Range("I54").Value = $ws.Range("I54").Value2 - 2
$ws.Range("I55").Value = 2.96
$ws.Range("I56").Value = $ws.Range("I56").Value2 - 2
$ws.Range("I57").Value = $ws.Range("I57").Value2 - 2
$ws.Range("I58").Value = 2.71
$ws.Range("I59").Value = $ws.Range("I59").Value2 - 2
$ws.Range("I60").Value = $ws.Range("I60").Value2 - 2
$ws.Range("I61").Value = 2.9800000000000004
$ws.Range("I62").Value = $ws.Range("I62").Value2 - 2
$ws.Range("I63").Value = $ws.Range("I63").Value2 - 2
$ws.Range("I64").Value = $ws.Range("I64").Value2 - 2
$ws.Range("I65").Value = 2.8599999999999994
$ws.Range("I66").Value = $ws.Range("I66").Value2 - 2
$ws.Range("I67").Value = $ws.Range("I67").Value2 - 2
$ws.Range("I68").Value = $ws.Range("I68").Value2 - 2
$ws.Range("I69").Value = $ws.Range("I69").Value2 - 2
$ws.Range("I70").Value = 2.9000000000000004
$ws.Range("I71").Value = $ws.Range("I71").Value2 - 2
$ws.Range("I72").Value = $ws.Range("I72").Value2 - 2
$ws.Range("I73").Value = $ws.Range("I73").Value2 - 2
$ws.Range("I74").Value = $ws.Range("I74").Value2 - 2
$ws.Range("I75").Value = $ws.Range("I75").Value2 - 2
$ws.Range("I76").Value = $ws.Range("I76").Value2 - 2
$ws.Range("I77").Value = $ws.Range("I77").Value2 - 2
$ws.Range("I78").Value = $ws.Range("I78").Value2 - 2
$ws.Range("I79").Value = $ws.Range("I79").Value2 - 2
$ws.Range("I80").Value = 2.8100000000000005
$ws.Range("I81").Value = $ws.Range("I81").Value2 - 2
$ws.Range("I82").Value = $ws.Range("I82").Value2 - 2
$ws.Range("I83").Value = $ws.Range("I83").Value2 - 2
$ws.Range("I84").Value = $ws.Range("I84").Value2 - 2
$ws.Range("I85").Value = $ws.Range("I85").Value2 - 2
$ws.Range("I86").Value = $ws.Range("I86").Value2 - 2
$ws.Range("I87").Value = $ws.Range("I87").Value2 - 2
$ws.Range("I88").Value = $ws.Range("I88").Value2 - 2
$ws.Range("I89").Value = $ws.Range("I89").Value2 - 2
$ws.Range("I90").Value = 2.8
$ws.Range("I91").Value = 2.7
$ws.Range("I92").Value = $ws.Range("I92").Value2 - 2
$ws.Range("I93").Value = $ws.Range("I93").Value2 - 2
$ws.Range("I94").Value = $ws.Range("I94").Value2 - 2
$ws.Range("I95").Value = $ws.Range("I95").Value2 - 2
$ws.Range("I96").Value = $ws.Range("I96").Value2 - 2
$ws.Range("I97").Value = $ws.Range("I97").Value2 - 2
$ws.Range("I98").Value = 2.71
$ws.Range("I99").Value = $ws.Range("I99").Value2 - 2
$ws.Range("I100").Value = $ws.Range("I100").Value2 - 2
$ws.Range("I101").Value = $ws.Range("I101").Value2 - 2
$ws.Range("I102").Value = 2.7
$ws.Range("I103").Value = $ws.Range("I103").Value2 - 2
$ws.Range("I104").Value = 2.7
$ws.Range("I105").Value = $ws.Range("I105").Value2 - 2
$ws.Range("I106").Value = $ws.Range("I106").Value2 - 2
$ws.Range("I107").Value = $ws.Range("I107").Value2 - 2
$ws.Range("I108").Value = $ws.Range("I108").Value2 - 2
$ws.Range("I109").Value = $ws.Range("I109").Value2 - 2
$ws.Range("I110").Value = $ws.Range("I110").Value2 - 2
$ws.Range("I111").Value = 2.9299999999999997
$ws.Range("I112").Value = $ws.Range("I112").Value2 - 2
